# Add 2022-Q4 data:
#  1. Insert a new "2022-Q4" worksheet right after "总计", built from a
#     duplicate of "2022-Q2" (to inherit identical column layout / styles),
#     then overwrite its single data row with the 2022-Q4 fund holding.
#  2. Insert a new row into "总计" for the 2022-Q4 summary figures.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2022-Q2")

# --- 1. Build the new "2022-Q4" sheet -------------------------------------
$templateSheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Template has 4 data rows (2..5); 2022-Q4 only needs one (row 2).
$q4Sheet.Rows("3:5").Delete()

$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'015729"
$q4Sheet.Range("C2").Value = "朱雀碳中和三年持有期混合"
$q4Sheet.Range("D2").Value = "'3.60"
$q4Sheet.Range("E2").Value = "'49.87"
$q4Sheet.Range("F2").Value = "'2.28"
$q4Sheet.Range("G2").Value = "'0.0821"
$q4Sheet.Range("H2").Value = 8
# Drop the "number stored as text" quote-prefix marker picked up above while
# keeping the values themselves as text.
$q4Sheet.Range("B2:G2").ClearFormats()

# --- 2. Insert the summary row into "总计" ---------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.08

$totalSheet.Activate()
